$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 142.9073533333333
$ws.Range("H2").Value = 428.72206
$ws.Range("I2").Value = 0.5576664151504187
$ws.Range("J2").Value = 0.5576664151504188
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 103.4766596666667
$ws.Range("N2").Value = 310.429979
$ws.Range("O2").Value = 0.877785331764719
$ws.Range("P2").Value = 0.8777853317647188
$ws.Range("Q2").Value = 14787.57556473742
$ws.Range("R2").Value = 133088.1800826367
$ws.Range("S2").Value = 0.4895113992368518
$ws.Range("T2").Value = 0.4895113992368518

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 142.9073533333333
$ws.Range("H3").Value = 428.72206
$ws.Range("I3").Value = 0.5576664151504187
$ws.Range("J3").Value = 0.5576664151504188
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8265796666666668
$ws.Range("N3").Value = 2.479739
$ws.Range("O3").Value = 0.007011818020336602
$ws.Range("P3").Value = 0.0070118180203366
$ws.Range("Q3").Value = 118.1243124824822
$ws.Range("R3").Value = 1063.11881234234
$ws.Range("S3").Value = 0.003910255419088219
$ws.Range("T3").Value = 0.003910255419088218

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 142.9073533333333
$ws.Range("H4").Value = 428.72206
$ws.Range("I4").Value = 0.5576664151504187
$ws.Range("J4").Value = 0.5576664151504188
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.58054833333333
$ws.Range("N4").Value = 40.741645
$ws.Range("O4").Value = 0.1152028502149446
$ws.Range("P4").Value = 0.1152028502149446
$ws.Range("Q4").Value = 1940.760219132078
$ws.Range("R4").Value = 17466.8419721887
$ws.Range("S4").Value = 0.0642447604944788
$ws.Range("T4").Value = 0.0642447604944788

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 63.967809
$ws.Range("H5").Value = 191.903427
$ws.Range("I5").Value = 0.2496211559306514
$ws.Range("J5").Value = 0.2496211559306514
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 103.4766596666667
$ws.Range("N5").Value = 310.429979
$ws.Range("O5").Value = 0.877785331764719
$ws.Range("P5").Value = 0.8777853317647188
$ws.Range("Q5").Value = 6619.175201515336
$ws.Range("R5").Value = 59572.57681363803
$ws.Range("S5").Value = 0.2191137891740794
$ws.Range("T5").Value = 0.2191137891740794

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 63.967809
$ws.Range("H6").Value = 191.903427
$ws.Range("I6").Value = 0.2496211559306514
$ws.Range("J6").Value = 0.2496211559306514
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8265796666666668
$ws.Range("N6").Value = 2.479739
$ws.Range("O6").Value = 0.007011818020336602
$ws.Range("P6").Value = 0.0070118180203366
$ws.Range("Q6").Value = 52.87449024061701
$ws.Range("R6").Value = 475.870412165553
$ws.Range("S6").Value = 0.001750298119411794
$ws.Range("T6").Value = 0.001750298119411794

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 63.967809
$ws.Range("H7").Value = 191.903427
$ws.Range("I7").Value = 0.2496211559306514
$ws.Range("J7").Value = 0.2496211559306514
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.58054833333333
$ws.Range("N7").Value = 40.741645
$ws.Range("O7").Value = 0.1152028502149446
$ws.Range("P7").Value = 0.1152028502149446
$ws.Range("Q7").Value = 868.7179219019349
$ws.Range("R7").Value = 7818.461297117415
$ws.Range("S7").Value = 0.02875706863716016
$ws.Range("T7").Value = 0.02875706863716016

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 49.38440333333333
$ws.Range("H8").Value = 148.15321
$ws.Range("I8").Value = 0.1927124289189298
$ws.Range("J8").Value = 0.1927124289189298
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.4766596666667
$ws.Range("N8").Value = 310.429979
$ws.Range("O8").Value = 0.877785331764719
$ws.Range("P8").Value = 0.8777853317647188
$ws.Range("Q8").Value = 5110.133096564732
$ws.Range("R8").Value = 45991.19786908259
$ws.Range("S8").Value = 0.1691601433537876
$ws.Range("T8").Value = 0.1691601433537876

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 49.38440333333333
$ws.Range("H9").Value = 148.15321
$ws.Range("I9").Value = 0.1927124289189298
$ws.Range("J9").Value = 0.1927124289189298
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8265796666666668
$ws.Range("N9").Value = 2.479739
$ws.Range("O9").Value = 0.007011818020336602
$ws.Range("P9").Value = 0.0070118180203366
$ws.Range("Q9").Value = 40.82014364579889
$ws.Range("R9").Value = 367.3812928121901
$ws.Range("S9").Value = 0.001351264481836588
$ws.Range("T9").Value = 0.001351264481836588

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 49.38440333333333
$ws.Range("H10").Value = 148.15321
$ws.Range("I10").Value = 0.1927124289189298
$ws.Range("J10").Value = 0.1927124289189298
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.58054833333333
$ws.Range("N10").Value = 40.741645
$ws.Range("O10").Value = 0.1152028502149446
$ws.Range("P10").Value = 0.1152028502149446
$ws.Range("Q10").Value = 670.6672763811611
$ws.Range("R10").Value = 6036.00548743045
$ws.Range("S10").Value = 0.02220102108330563
$ws.Range("T10").Value = 0.02220102108330563
